$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns so the "4_train" group gets its own GNN-MT columns,
# matching the layout already used by the 8_train / 16_train groups.
$ws.Columns("C").Insert()
$ws.Columns("F").Insert()

# Header row: label the two newly inserted columns
$ws.Range("C1").Value = "4_train (GNN-MT)"
$ws.Range("F1").Value = "4_train (GNN-MT) std"

# Refresh the results table (rows 2-5) with the updated aggregate values
# Row 2: Viruses
$ws.Range("B2").Value = 0.2221029411764706
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0.1127058823529411
$ws.Range("E2").Value = 0.0147381623389696
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0.0140127786282842
$ws.Range("H2").Value = 0.2077258064516129
$ws.Range("I2").Value = 0.1797740740740741
$ws.Range("J2").Value = 0.156074193548387
$ws.Range("K2").Value = 0.0148987010327458
$ws.Range("L2").Value = 0.0135619692803948
$ws.Range("M2").Value = 0.013370374111326
$ws.Range("N2").Value = 0.192171568627451
$ws.Range("O2").Value = 0.1661186274509803
$ws.Range("P2").Value = 0.1856892156862744
$ws.Range("Q2").Value = 0.0118697072972182
$ws.Range("R2").Value = 0.0107955207220888
$ws.Range("S2").Value = 0.0122248428290348

# Row 3: Bacteria
$ws.Range("B3").Value = 0.2128
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0.0916857142857142
$ws.Range("E3").Value = 0.0340240350044375
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0.0293686442129214
$ws.Range("H3").Value = 0.2160428571428571
$ws.Range("I3").Value = 0.2126166666666667
$ws.Range("J3").Value = 0.1737857142857143
$ws.Range("K3").Value = 0.0335000796735657
$ws.Range("L3").Value = 0.0304324106632048
$ws.Range("M3").Value = 0.0314539142156869
$ws.Range("N3").Value = 0.1990873015873016
$ws.Range("O3").Value = 0.1840396825396825
$ws.Range("P3").Value = 0.2041825396825397
$ws.Range("Q3").Value = 0.0346461618407768
$ws.Range("R3").Value = 0.0267363020725282
$ws.Range("S3").Value = 0.0343227111299082

# Row 4: Fungi
$ws.Range("B4").Value = 0.334875
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0.265375
$ws.Range("E4").Value = 0.0428291612522522
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0.0669050333094817
$ws.Range("H4").Value = 0.1733125
$ws.Range("I4").Value = 0.1665
$ws.Range("J4").Value = 0.1833125
$ws.Range("K4").Value = 0.0819327991035849
$ws.Range("L4").Value = 0.0220510140860283
$ws.Range("M4").Value = 0.0969847570124188
$ws.Range("N4").Value = 0.2342708333333333
$ws.Range("O4").Value = 0.2337708333333333
$ws.Range("P4").Value = 0.2250208333333333
$ws.Range("Q4").Value = 0.0589741427076464
$ws.Range("R4").Value = 0.0759957842931045
$ws.Range("S4").Value = 0.1065382113625157

# Row 5: all
$ws.Range("B5").Value = 0.2253240223463687
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0.1154189944134078
$ws.Range("E5").Value = 0.0132233095837899
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0.0126303505116116
$ws.Range("H5").Value = 0.2078055555555555
$ws.Range("I5").Value = 0.1851794117647058
$ws.Range("J5").Value = 0.1603055555555555
$ws.Range("K5").Value = 0.0134170572563961
$ws.Range("L5").Value = 0.012048995224503
$ws.Range("M5").Value = 0.0123865917732178
$ws.Range("N5").Value = 0.1937799145299145
$ws.Range("O5").Value = 0.1694363247863248
$ws.Range("P5").Value = 0.1884876068376068
$ws.Range("Q5").Value = 0.0110285750778723
$ws.Range("R5").Value = 0.0099439697445949
$ws.Range("S5").Value = 0.0114232188759778

# Restore the last active cell selection recorded in the workbook
$ws.Range("I17").Select() | Out-Null
